$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the J column (constant 0.2 across all rows)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# New bold-12pt, vertically-centered style used only for the B-column values
$valueRange = $ws.Range("B14:B17")
$valueRange.Font.Bold = $true
$valueRange.Font.Size = 12
$valueRange.VerticalAlignment = -4108  # xlCenter

# Selection now covers the new summary block
$ws.Range("A14:B17").Select()

# Page setup matching the updated worksheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
